$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 with same formatting (bold, border, centered) as the
# other header cells by copying the format from E1, then set its value/text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate F2:F62 with the recorded time_taken timestamp values (plain text).
$timeTaken = @(
    "2021-10-05 13:42:33.641524",
    "2021-10-05 13:42:33.641537",
    "2021-10-05 13:42:33.641541",
    "2021-10-05 13:42:33.641545",
    "2021-10-05 13:42:33.641548",
    "2021-10-05 13:42:33.641551",
    "2021-10-05 13:42:33.641554",
    "2021-10-05 13:42:33.641557",
    "2021-10-05 13:42:33.641560",
    "2021-10-05 13:42:33.641564",
    "2021-10-05 13:42:33.641566",
    "2021-10-05 13:42:33.641569",
    "2021-10-05 13:42:33.641572",
    "2021-10-05 13:42:33.641575",
    "2021-10-05 13:42:33.641578",
    "2021-10-05 13:42:33.641581",
    "2021-10-05 13:42:33.641584",
    "2021-10-05 13:42:33.641587",
    "2021-10-05 13:42:33.641590",
    "2021-10-05 13:42:33.641593",
    "2021-10-05 13:42:33.641596",
    "2021-10-05 13:42:33.641599",
    "2021-10-05 13:42:33.641602",
    "2021-10-05 13:42:33.641605",
    "2021-10-05 13:42:33.641608",
    "2021-10-05 13:42:33.641611",
    "2021-10-05 13:42:33.641614",
    "2021-10-05 13:42:33.641617",
    "2021-10-05 13:42:33.641620",
    "2021-10-05 13:42:33.641623",
    "2021-10-05 13:42:33.641626",
    "2021-10-05 13:42:33.641629",
    "2021-10-05 13:42:33.641633",
    "2021-10-05 13:42:33.641636",
    "2021-10-05 13:42:33.641639",
    "2021-10-05 13:42:33.641642",
    "2021-10-05 13:42:33.641645",
    "2021-10-05 13:42:33.641648",
    "2021-10-05 13:42:33.641651",
    "2021-10-05 13:42:33.641654",
    "2021-10-05 13:42:33.641657",
    "2021-10-05 13:42:33.641660",
    "2021-10-05 13:42:33.641663",
    "2021-10-05 13:42:33.641666",
    "2021-10-05 13:42:33.641669",
    "2021-10-05 13:42:33.641672",
    "2021-10-05 13:42:33.641675",
    "2021-10-05 13:42:33.641678",
    "2021-10-05 13:42:33.641682",
    "2021-10-05 13:42:33.641684",
    "2021-10-05 13:42:33.641688",
    "2021-10-05 13:42:33.641691",
    "2021-10-05 13:42:33.641694",
    "2021-10-05 13:42:33.641697",
    "2021-10-05 13:42:33.641700",
    "2021-10-05 13:42:33.641703",
    "2021-10-05 13:42:33.641706",
    "2021-10-05 13:42:33.641709",
    "2021-10-05 13:42:33.641712",
    "2021-10-05 13:42:33.641715",
    "2021-10-05 13:42:33.641718"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

Write-Host "Added time_taken column F1:F62"
